$d = $word.ActiveDocument

# 1. Change "2 minutter" -> "4 minutter" (startup timer)
$d.Content.Find.Execute("2 minutter", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4 minutter", 2)

# 2. Swap 20 sekund <-> 15 sek
$d.Content.Find.Execute("hvert 20 sekund", $true, $false, $false, $false, $false,
                         $true, 1, $false, "hvert 15 sekund", 2)
$d.Content.Find.Execute("e.g. 15 sek", $true, $false, $false, $false, $false,
                         $true, 1, $false, "e.g. 20 sek", 2)
